# Add rare weapons to level 3
# - item-chances: two new loot-table rows for the rare sword / rare dagger (level 3)
# - stuff-descriptor: two new item rows describing the rare sword / rare dagger
# - restore the view/selection state recorded for each sheet

$wb = $excel.ActiveWorkbook

$wsEnemy = $wb.Worksheets.Item("enemy-chances")
$wsItem  = $wb.Worksheets.Item("item-chances")
$wsStuff = $wb.Worksheets.Item("stuff-descriptor")

# ---------------------------------------------------------------------------
# item-chances: new level-3 loot entries for RareSword / RareDagger
# ---------------------------------------------------------------------------
$wsItem.Range("A11").Value = 3
$wsItem.Range("B11").Value = "RareSword"
$wsItem.Range("C11").Value = 10

$wsItem.Range("A12").Value = 3
$wsItem.Range("B12").Value = "RareDagger"
$wsItem.Range("C12").Value = 15

# ---------------------------------------------------------------------------
# stuff-descriptor: new item rows for RareSword / RareDagger
# ---------------------------------------------------------------------------
$wsStuff.Range("A22").Value = "RareSword"
$wsStuff.Range("B22").Value = "sword"
$wsStuff.Range("C22").Value = "Rare Sword"
$wsStuff.Range("D22").Formula = '="Rare sword dealing "&H22&" damage"'
$wsStuff.Range("E22").Value = "#0d09ed"
$wsStuff.Range("H22").Value = 7
$wsStuff.Range("I22").Value = 0
$wsStuff.Range("I22").Style = $wsStuff.Range("I11").Style

$wsStuff.Range("A23").Value = "RareDagger"
$wsStuff.Range("B23").Value = "dagger"
$wsStuff.Range("C23").Value = "Rare Dagger"
$wsStuff.Range("D23").Formula = '="Rare dagger dealing "&H23&" damage"'
$wsStuff.Range("E23").Value = "#00BFFF"
$wsStuff.Range("H23").Formula = '=FLOOR.MATH((H22+H11)/2)'
$wsStuff.Range("I23").Value = 1
$wsStuff.Range("I23").Style = $wsStuff.Range("I11").Style

# trailing, formatted-but-empty cells left behind at the bottom of the table
$wsStuff.Range("J24").Style = $wsStuff.Range("J9").Style
$wsStuff.Range("J25").Style = $wsStuff.Range("J9").Style

# ---------------------------------------------------------------------------
# Restore per-sheet selection / view state
# ---------------------------------------------------------------------------
$null = $wsEnemy.Range("A10").Select()
$null = $wsStuff.Range("A24").Select()

$null = $wsItem.Range("C12").Select()
$null = $wsItem.Activate()

Write-Output "done"
